# Add new register rows documenting the gpio monitor/control registers
# (gpiocontrol, gpio_in_count, gpio_out_count) to the cr_registers sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column A (registername) for the new gpiocontrol rows -------------
$ws.Cells.Item(40, 1).Value = "gpio_in_enable"
$ws.Cells.Item(41, 1).Value = "gpio_out_enable"
$ws.Cells.Item(42, 1).Value = "gpio_in_count_rst"
$ws.Cells.Item(43, 1).Value = "gpio_out_count_rst"

# --- column G (description) for all six new rows -----------------------
$ws.Cells.Item(40, 7).Value = "Value 1 enables using the gpio_in pin"
$ws.Cells.Item(41, 7).Value = "Value 1 enables using the gpio_out pin"
$ws.Cells.Item(42, 7).Value = "Changing the value from 0 to 1 resets the gpio in counter."
$ws.Cells.Item(43, 7).Value = "Changing the value from 0 to 1 resets the gpio out counter."
$ws.Cells.Item(44, 7).Value = "Counts the number of rising edges on the gpio in pin, even if the pin is disable from sending a signal to the rest of the firmware."
$ws.Cells.Item(45, 7).Value = "Counts the number of rising edges on the gpio out pin, even if the pin is disable from sending a signal to the rest of the firmware."

# --- column A for the gpio_in_count row --------------------------------
$ws.Cells.Item(44, 1).Value = "gpio_in_count"

# --- column C (mainregister) for the gpiocontrol rows -------------------
$ws.Cells.Item(40, 3).Value = "cosmic_ray_gpiocontrol"
$ws.Cells.Item(41, 3).Value = "cosmic_ray_gpiocontrol"
$ws.Cells.Item(42, 3).Value = "cosmic_ray_gpiocontrol"
$ws.Cells.Item(43, 3).Value = "cosmic_ray_gpiocontrol"

# --- column C for the gpio_in_count row ---------------------------------
$ws.Cells.Item(44, 3).Value = "cosmic_ray_gpio_in_count"

# --- column A for the gpio_out_count row --------------------------------
$ws.Cells.Item(45, 1).Value = "gpio_out_count"

# --- column C for the gpio_out_count row --------------------------------
$ws.Cells.Item(45, 3).Value = "cosmic_ray_gpio_out_count"

# --- remaining numeric / reused-string columns --------------------------
# bitwidth (B), offset_from_msb (D), mainregister_bitwidth (E), readonly/readwrite (F)
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 32
$ws.Cells.Item(40, 6).Value = "readwrite"

$ws.Cells.Item(41, 2).Value = 1
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 5).Value = 32
$ws.Cells.Item(41, 6).Value = "readwrite"

$ws.Cells.Item(42, 2).Value = 1
$ws.Cells.Item(42, 4).Value = 2
$ws.Cells.Item(42, 5).Value = 32
$ws.Cells.Item(42, 6).Value = "readwrite"

$ws.Cells.Item(43, 2).Value = 1
$ws.Cells.Item(43, 4).Value = 3
$ws.Cells.Item(43, 5).Value = 32
$ws.Cells.Item(43, 6).Value = "readwrite"

$ws.Cells.Item(44, 2).Value = 32
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 32
$ws.Cells.Item(44, 6).Value = "readonly"

$ws.Cells.Item(45, 2).Value = 32
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 32
$ws.Cells.Item(45, 6).Value = "readonly"

# --- update the view: scroll/selection to reflect the new rows ----------
$ws.Range("F47").Select() | Out-Null
